$wb = $excel.ActiveWorkbook

# Update the "想去人数" (want-to-go count) figures on the "展览" sheet
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1490
$ws1.Range("F3").Value = 3141
$ws1.Range("F5").Value = 876

# Same updates are mirrored on the "全部类型" sheet
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1490
$ws4.Range("F3").Value = 3141
$ws4.Range("F5").Value = 876
